$d = $word.ActiveDocument

# The document contains two "EDU_DATE" placeholders (split across an "E"
# run and a "DU_DATE" run each time):
#   1) "...기  간   :   EDU_DATE"                 (left untouched)
#   2) "...주관하는EDU_DATE 교육 및 실습과정을..."  (renamed to COURSE)
#
# Locate the first occurrence purely to know where it ends, so the next
# Find can be restricted to the remainder of the document and therefore
# only ever hit the second occurrence.
$r1 = $d.Content
$found1 = $r1.Find.Execute("EDU_DATE", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)

if (-not $found1) {
    throw "Could not find first EDU_DATE occurrence"
}

$r2 = $d.Range($r1.End, $d.Content.End)
$found2 = $r2.Find.Execute("EDU_DATE", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)

if (-not $found2) {
    throw "Could not find second EDU_DATE occurrence"
}

# r2 now spans the second "EDU_DATE" (an "E" run followed by a "DU_DATE"
# run). Replace the text of the "DU_DATE" run with "COURSE" (this keeps
# that run's original formatting -- no eastAsia font hint, just sz/szCs),
# then delete the now-redundant leading "E" run.
$duDateRange = $d.Range($r2.Start + 1, $r2.End)
$duDateRange.Text = "COURSE"

$eRange = $d.Range($r2.Start, $r2.Start + 1)
$eRange.Delete()
